# Commit: "update file with jgit"
# The "Rules" sheet's E8 cell held the greeting text for rule R10
# ("Good Morning"). The commit replaces that text with "GIT UPDATE".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"

# Leave the edited cell selected, matching the saved view state.
$ws.Range("E8").Select()
